{"js": "// \"mais um filme adicionado\" - add a new movie (\"Entre facas\") to the list,\n// right after \"Poderoso Chef\u00e3o \" and before the trailing empty paragraph.\n//\n// Touching the first two paragraphs (\"Monkey Man\" / \"Past Lives \") also\n// makes Word clean up their stale spell-check (proofErr) markers and merge\n// their split runs into a single run, which is reflected in the target\n// document too.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Normalize \"Monkey\" + \" Man\" -> \"Monkey Man\" (clears stray proofErr marks\n// and merges the two runs into one).\nconst monkeyPara = paragraphs.items.find(\n  (p) => p.text.indexOf(\"Monkey\") !== -1\n);\nif (monkeyPara) {\n  monkeyPara.clear();\n  monkeyPara.insertText(\"Monkey Man\", \"Start\");\n}\n\n// Normalize \"Past\" + \" Lives \" -> \"Past Lives \" (same idea).\nconst pastPara = paragraphs.items.find((p) => p.text.indexOf(\"Past\") !== -1);\nif (pastPara) {\n  pastPara.clear();\n  pastPara.insertText(\"Past Lives \", \"Start\");\n}\n\nawait context.sync();\n\n// Find the paragraph that currently ends the movie list (\"Poderoso Chef\u00e3o \")\n// and add the new movie right after it, before the trailing blank paragraph.\nconst chefaoParagraph = paragraphs.items.find(\n  (p) => p.text.indexOf(\"Poderoso Chef\") !== -1\n);\n\nif (chefaoParagraph) {\n  chefaoParagraph.insertParagraph(\"Entre facas\", \"After\");\n} else {\n  // Fallback: insert before the final (blank) paragraph.\n  const last = paragraphs.items[paragraphs.items.length - 1];\n  last.insertParagraph(\"Entre facas\", \"Before\");\n}\n\nawait context.sync();\n", "ps1": "# \"mais um filme adicionado\" - add a new movie (\"Entre facas\") to the list,\n# right after \"Poderoso Chefao \" and before the trailing empty paragraph.\n#\n# Touching the first two paragraphs (\"Monkey Man\" / \"Past Lives \") also\n# makes Word clean up their stale spell-check (proofErr) markers and merge\n# their split runs into a single run, which is reflected in the target\n# document too. Deleting the whole paragraph range (including its paragraph\n# mark) and re-inserting a fresh paragraph is what actually clears those\n# stray proofErr markers (simply overwriting .Text leaves them behind).\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($needle) {\n  for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"$needle*\") {\n      return $i\n    }\n  }\n  return -1\n}\n\n# --- Normalize \"Monkey\" + \" Man\" -> \"Monkey Man\" ---\n$monkeyIndex = Find-ParagraphIndex \"Monkey\"\nif ($monkeyIndex -ne -1) {\n  $d.Paragraphs.Item($monkeyIndex).Range.Delete()\n  $d.Paragraphs.Item($monkeyIndex).Range.InsertParagraphBefore()\n  $d.Paragraphs.Item($monkeyIndex).Range.Text = \"Monkey Man\"\n}\n\n# --- Normalize \"Past\" + \" Lives \" -> \"Past Lives \" ---\n$pastIndex = Find-ParagraphIndex \"Past\"\nif ($pastIndex -ne -1) {\n  $d.Paragraphs.Item($pastIndex).Range.Delete()\n  $d.Paragraphs.Item($pastIndex).Range.InsertParagraphBefore()\n  $d.Paragraphs.Item($pastIndex).Range.Text = \"Past Lives \"\n}\n\n# --- Find the paragraph that ends the movie list (\"Poderoso Chefao \") and\n#     insert the new movie right after it, before the trailing blank\n#     paragraph. ---\n$targetIndex = Find-ParagraphIndex \"Poderoso Chef\"\n\nif ($targetIndex -eq -1) {\n  # Fallback: insert before the final (blank) paragraph.\n  $targetIndex = $d.Paragraphs.Count - 1\n}\n\n$targetRange = $d.Paragraphs.Item($targetIndex).Range\n$targetRange.InsertParagraphAfter()\n$d.Paragraphs.Item($targetIndex + 1).Range.Text = \"Entre facas\"\n"}
